$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old rows (4 through 18) that are being removed
$ws.Range("A4:A18").ClearContents()

# Update remaining cells:
# A1: "Plots anzeigen" -> "Plots anzeigen/ausblenden"
$ws.Range("A1").Value = "Plots anzeigen/ausblenden"
# A2: unchanged - "Umschaltung AFR <> Lambda"
# A3: now holds what used to be in A4 - "Cursor-Werte"
$ws.Range("A3").Value = "Cursor-Werte"

# Selection moves to A3
$ws.Range("A3").Select()
